$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value in a cell guaranteed to be stored as a shared string
# (avoids Excel auto-converting numeric-looking text, and avoids creating new
# cell styles / number formats in the process).
function Set-TextValue($ws, $cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.EntireColumn.Delete()
}

# Source row to copy formatting/styles from for each new row
$srcRow = 323

foreach ($r in 324..335) {
    $ws.Range("A$srcRow" + ":M$srcRow").Copy($ws.Range("A$r" + ":M$r"))
    $ws.Rows.Item($r).RowHeight = 16
}

# Row 324
Set-TextValue $ws "A324" '●'
Set-TextValue $ws "D324" '10703'
Set-TextValue $ws "E324" 'Event month'
Set-TextValue $ws "F324" '1: 1414'
Set-TextValue $ws "G324" '1: 1421'
$ws.Range("H324").Value = 0
Set-TextValue $ws "I324" 'December'
$ws.Range("J324").Value = 8
$ws.Range("K324").Value = 0.016161
Set-TextValue $ws "L324" 'Sonia'
Set-TextValue $ws "M324" '11/14/18 10:06:00'

# Row 325
Set-TextValue $ws "A325" '●'
Set-TextValue $ws "D325" '10703'
Set-TextValue $ws "E325" 'Event month'
Set-TextValue $ws "F325" '1: 1431'
Set-TextValue $ws "G325" '1: 1433'
$ws.Range("H325").Value = 0
Set-TextValue $ws "I325" 'May'
$ws.Range("J325").Value = 3
$ws.Range("K325").Value = 0.00606
Set-TextValue $ws "L325" 'Sonia'
Set-TextValue $ws "M325" '11/14/18 10:06:00'

# Row 326
Set-TextValue $ws "A326" '●'
Set-TextValue $ws "D326" '10703'
Set-TextValue $ws "E326" 'Event year'
Set-TextValue $ws "F326" '1: 1423'
Set-TextValue $ws "G326" '1: 1426'
$ws.Range("H326").Value = 0
Set-TextValue $ws "I326" '2004'
$ws.Range("J326").Value = 4
$ws.Range("K326").Value = 0.008081
Set-TextValue $ws "L326" 'Sonia'
Set-TextValue $ws "M326" '11/14/18 10:06:00'

# Row 327
Set-TextValue $ws "A327" '●'
Set-TextValue $ws "D327" '10703'
Set-TextValue $ws "E327" 'Event year'
Set-TextValue $ws "F327" '1: 1435'
Set-TextValue $ws "G327" '1: 1438'
$ws.Range("H327").Value = 0
Set-TextValue $ws "I327" '2005'
$ws.Range("J327").Value = 4
$ws.Range("K327").Value = 0.008081
Set-TextValue $ws "L327" 'Sonia'
Set-TextValue $ws "M327" '11/14/18 10:06:00'

# Row 328
Set-TextValue $ws "A328" '●'
Set-TextValue $ws "D328" '10703'
Set-TextValue $ws "E328" 'A'
Set-TextValue $ws "F328" '1: 1414'
Set-TextValue $ws "G328" '1: 1421'
$ws.Range("H328").Value = 0
Set-TextValue $ws "I328" 'December'
$ws.Range("J328").Value = 8
$ws.Range("K328").Value = 0.016161
Set-TextValue $ws "L328" 'Sonia'
Set-TextValue $ws "M328" '11/14/18 10:06:00'

# Row 329
Set-TextValue $ws "A329" '●'
Set-TextValue $ws "D329" '10703'
Set-TextValue $ws "E329" 'A'
Set-TextValue $ws "F329" '1: 1423'
Set-TextValue $ws "G329" '1: 1426'
$ws.Range("H329").Value = 0
Set-TextValue $ws "I329" '2004'
$ws.Range("J329").Value = 4
$ws.Range("K329").Value = 0.008081
Set-TextValue $ws "L329" 'Sonia'
Set-TextValue $ws "M329" '11/14/18 10:06:00'

# Row 330
Set-TextValue $ws "A330" '●'
Set-TextValue $ws "D330" '10703'
Set-TextValue $ws "E330" 'B'
Set-TextValue $ws "F330" '1: 1431'
Set-TextValue $ws "G330" '1: 1433'
$ws.Range("H330").Value = 0
Set-TextValue $ws "I330" 'May'
$ws.Range("J330").Value = 3
$ws.Range("K330").Value = 0.00606
Set-TextValue $ws "L330" 'Sonia'
Set-TextValue $ws "M330" '11/14/18 10:07:00'

# Row 331
Set-TextValue $ws "A331" '●'
Set-TextValue $ws "D331" '10703'
Set-TextValue $ws "E331" 'B'
Set-TextValue $ws "F331" '1: 1435'
Set-TextValue $ws "G331" '1: 1438'
$ws.Range("H331").Value = 0
Set-TextValue $ws "I331" '2005'
$ws.Range("J331").Value = 4
$ws.Range("K331").Value = 0.008081
Set-TextValue $ws "L331" 'Sonia'
Set-TextValue $ws "M331" '11/14/18 10:07:00'

# Row 332
Set-TextValue $ws "A332" '●'
Set-TextValue $ws "D332" '11202'
Set-TextValue $ws "E332" 'Event month'
Set-TextValue $ws "F332" '1: 506'
Set-TextValue $ws "G332" '1: 513'
$ws.Range("H332").Value = 0
Set-TextValue $ws "I332" 'December'
$ws.Range("J332").Value = 8
$ws.Range("K332").Value = 0.041216
Set-TextValue $ws "L332" 'Sonia'
Set-TextValue $ws "M332" '11/14/18 10:07:00'

# Row 333
Set-TextValue $ws "A333" '●'
Set-TextValue $ws "D333" '11202'
Set-TextValue $ws "E333" 'Event month'
Set-TextValue $ws "F333" '1: 523'
Set-TextValue $ws "G333" '1: 530'
$ws.Range("H333").Value = 0
Set-TextValue $ws "I333" 'December'
$ws.Range("J333").Value = 8
$ws.Range("K333").Value = 0.041216
Set-TextValue $ws "L333" 'Sonia'
Set-TextValue $ws "M333" '11/14/18 10:07:00'

# Row 334
Set-TextValue $ws "A334" '●'
Set-TextValue $ws "D334" '11202'
Set-TextValue $ws "E334" 'Event year'
Set-TextValue $ws "F334" '1: 515'
Set-TextValue $ws "G334" '1: 518'
$ws.Range("H334").Value = 0
Set-TextValue $ws "I334" '2011'
$ws.Range("J334").Value = 4
$ws.Range("K334").Value = 0.020608
Set-TextValue $ws "L334" 'Sonia'
Set-TextValue $ws "M334" '11/14/18 10:07:00'

# Row 335
Set-TextValue $ws "A335" '●'
Set-TextValue $ws "D335" '11202'
Set-TextValue $ws "E335" 'Event year'
Set-TextValue $ws "F335" '1: 532'
Set-TextValue $ws "G335" '1: 535'
$ws.Range("H335").Value = 0
Set-TextValue $ws "I335" '2012'
$ws.Range("J335").Value = 4
$ws.Range("K335").Value = 0.020608
Set-TextValue $ws "L335" 'Sonia'
Set-TextValue $ws "M335" '11/14/18 10:08:00'
